# Correct mismatch in subscript mapping in BCS file
#
# The "BCS-BCS" sheet had the duration-multiplier formula
# (About!$B$11*About!$A$9*(1-Calculations!<col>3)) wired to row 2
# ("electricity sector") instead of row 3 ("industry sector"), and the
# plain formula (About!$B$11*About!$A$9) wired to row 3 instead of row 2.
# This swaps the two rows' contents back to the correct mapping.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("BCS-BCS")

# --- Row 2 ("electricity sector"): plain subsidy value, no duration factor ---
# B2:C2 stay 0 (unchanged), D2:P2 get the plain formula, Q2:AE2 go flat to 0.
$ws.Range("D2:P2").FormulaR1C1 = "=About!R11C2*About!R9C1"
$ws.Range("Q2:AE2").Value = 0

# Electricity-sector row no longer carries any special number formatting
# (it previously inherited the "General w/ explicit format" style from the
# old formula cells).
$ws.Range("B2:AE2").ClearFormats()

# --- Row 3 ("industry sector"): duration-adjusted subsidy value ---
# B3:C3 stay 0 (unchanged), D3:AE3 get the duration-adjusted formula.
$ws.Range("D3:AE3").FormulaR1C1 = "=About!R11C2*About!R9C1*(1-Calculations!RC)"

# Industry-sector row keeps no special formatting either.
$ws.Range("B3:AE3").ClearFormats()

# --- View state: BCS-BCS becomes the active sheet/tab, with a new selection ---
$ws.Activate()
$ws.Range("M2:P2").Select()
